# The data rows (2-14) were reshuffled (e.g. the sheet was re-sorted), so
# each row's content in columns D and K:S ends up equal to some other row's
# original content. Columns A,B,C,E,F,G,H,I,J,T are identical for every row
# and are left untouched.
#
# Strategy: snapshot the D2:S14 block with .Value2 (1-based 2D SAFEARRAY),
# build a new 1-based array with the rows permuted according to the mapping
# below, then write it back in one shot with .Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:S14")
$arr = $rng.Value2

# New-row -> old-row mapping (1 = sheet row 2, 2 = sheet row 3, ... 13 = sheet row 14)
$map = @{
    1  = 5
    2  = 4
    3  = 12
    4  = 13
    5  = 10
    6  = 6
    7  = 8
    8  = 2
    9  = 1
    10 = 7
    11 = 9
    12 = 11
    13 = 3
}

$rows = 13
$cols = 16

# NB: New-Object creates a 0-based array, while Range.Value2 returns a
# 1-based SAFEARRAY - offset the write side by one so values line up.
$new = New-Object 'object[,]' $rows, $cols

for ($r = 1; $r -le $rows; $r++) {
    $srcRow = $map[$r]
    for ($c = 1; $c -le $cols; $c++) {
        $new[$r - 1, $c - 1] = $arr[$srcRow, $c]
    }
}

$rng.Value = $new
